$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.257.91'
$ws.Range("E2").Value = '  +0.02%  '

$ws.Range("D3").Value = '2.486.07'
$ws.Range("E3").Value = '  -1.16%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = "'321.41"
$ws.Range("E5").Value = '  -0.85%  '

$ws.Range("D6").Value = "'107.70"
$ws.Range("E6").Value = '  +2.18%  '

$ws.Range("D7").Value = "'0.520"

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("E9").Value = '  -1.46%  '

$ws.Range("D10").Value = "'38.46"
$ws.Range("E10").Value = '  +4.50%  '

$ws.Range("D11").Value = "'0.0807"
$ws.Range("E11").Value = '  -1.52%  '

$ws.Range("E12").Value = '  -0.06%  '

$ws.Range("D13").Value = "'18.25"
$ws.Range("E13").Value = '  -0.80%  '

$ws.Range("D14").Value = "'7.09"
$ws.Range("E14").Value = '  -1.34%  '

$ws.Range("D15").Value = '2.877.00'
$ws.Range("E15").Value = '  -1.06%  '

$ws.Range("D16").Value = '2.480.95'
$ws.Range("E16").Value = '  -2.41%  '

$ws.Range("E17").Value = '  -0.31%  '

$ws.Range("D18").Value = '47.165.24'
$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("E19").Value = '  -0.91%  '

$ws.Range("D20").Value = "'6.61"
$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("E21").Value = '  -1.38%  '

$ws.Range("D22").Value = "'2.70"
$ws.Range("E22").Value = '  +12.76%  '

$ws.Range("D23").Value = "'70.25"
$ws.Range("E23").Value = '  -1.00%  '

$ws.Range("D24").Value = "'245.08"

$ws.Range("E25").Value = '  +0.15%  '

$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("D27").Value = "'25.62"
$ws.Range("E27").Value = '  -3.76%  '

$ws.Range("E28").Value = '  +3.16%  '

$ws.Range("D29").Value = "'9.95"
$ws.Range("E29").Value = '  -0.40%  '

$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = "'0.135"
$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = "'34.28"
$ws.Range("E31").Value = '  -2.41%  '

$ws.Range("D32").Value = "'49.48"
$ws.Range("E32").Value = '  -0.52%  '

$ws.Range("D33").Value = "'20.10"
$ws.Range("E33").Value = '  +1.42%  '

$ws.Range("D34").Value = "'5.31"
$ws.Range("E34").Value = '  -0.25%  '

$ws.Range("D35").Value = "'0.0776"
$ws.Range("E35").Value = '  -1.47%  '

$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = '  +0.14%  '

$ws.Range("E37").Value = '  -0.20%  '

$ws.Range("D38").Value = "'4.63"
$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("E39").Value = '  -2.07%  '

$ws.Range("D40").Value = "'22.61"
$ws.Range("E40").Value = '  +5.35%  '

$ws.Range("E41").Value = '  -0.72%  '

$ws.Range("E42").Value = '  -0.96%  '

$ws.Range("D43").Value = "'118.57"
$ws.Range("E43").Value = '  -4.13%  '

$ws.Range("D44").Value = "'0.0295"
$ws.Range("E44").Value = '  -0.81%  '

$ws.Range("D45").Value = '1.983.79'
$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("E46").Value = '  -0.61%  '

$ws.Range("D47").Value = "'1.99"
$ws.Range("E47").Value = '  -7.25%  '

$ws.Range("E48").Value = '  -1.11%  '

$ws.Range("E49").Value = '  -2.50%  '

$ws.Range("D50").Value = "'5.11"
$ws.Range("E50").Value = '  -5.02%  '

$ws.Range("D51").Value = "'56.56"
